# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (clrScheme "Office")      -- used by the Notes Master
#   ppt/theme/theme2.xml -> "Integral"     (clrScheme "Red Violet")  -- used by the Slide Master
#       (i.e. by every visible slide in the deck)
#
# The authored edit swaps the two themes' contents in place (the parts keep their
# names/relationships; only what's *in* them changes), so the deck's visible slides
# end up using the Office palette and the notes side ends up using the Red Violet
# ("Integral") palette.
#
# PowerPoint's object model doesn't expose raw theme-part XML for rewriting, but it
# does expose every theme color as a read/write RGB value through ThemeColorScheme,
# so we reproduce the net effect on the slide-facing theme (theme2.xml) by recoloring
# each of its 12 slots to the Office theme's values, in the documented
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation

# Any slide's ThemeColorScheme reflects the one deck-wide theme used by the slide
# master, so Slide 1 is as good an anchor as any.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
